$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Goal: change "...(which are tuples of Char, Char, Int). The first..."
# into "...(which are tuples of Int, Char, Char). The first...", split
# across 4 runs with a relocated "_GoBack" bookmark, matching Word's
# own behaviour when a user edits text in the middle of a run (the
# untouched leading part of the run keeps its original rsid, while the
# freshly (re)written parts end up rsid-less).
# ------------------------------------------------------------------

# 1. Locate the run span that must stay completely untouched (rsid must
#    survive): from the start of the run up to and including "tuples".
$rWholeRun = $d.Content
$rWholeRun.Find.Execute(" letters of the alphabet", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runStart = $rWholeRun.Start

$rTuplesEnd = $d.Content
$rTuplesEnd.Find.Execute("which are tuples", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$keepEnd = $rTuplesEnd.End

# 2. Temporarily mark that span with a harmless character formatting
#    toggle so the engine cannot silently re-merge/rewrite it when we
#    touch the rest of the run next to it.
$rKeep = $d.Range($runStart, $keepEnd)
$rKeep.Bold = 1

# 3. Perform the actual text change further along in the same run:
#    "Char, Char, Int" -> "Int, Char, Char".
$rReorder = $d.Content
$rReorder.Find.Execute("Char, Char, Int", $true, $false, $false, $false, $false, $true, 1, $false, "Int, Char, Char", 2)

# 4. Remove the temporary formatting marker again.
$rKeep2 = $d.Range($runStart, $keepEnd)
$rKeep2.Bold = 0

# 5. Split the remainder into the three additional runs shown in the
#    diff, using bookmark insert/delete as a pure (non text-editing)
#    run-splitting mechanism.

# 5a. Split point between "...tuples of Int, Char, Char" and ")." :
$rSplitEnd = $d.Content
$rSplitEnd.Find.Execute("tuples of Int, Char, Char", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos1 = $rSplitEnd.End
$bm1 = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("TempSplitMark", $bm1)
$d.Bookmarks("TempSplitMark").Delete()

# 5b. Real "_GoBack" bookmark goes between "tuples of Int, " and "Char, Char".
#     Adding a bookmark named "_GoBack" automatically relocates (removes)
#     any other bookmark of that name elsewhere in the document.
$rGoBack = $d.Content
$rGoBack.Find.Execute("tuples of Int, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos2 = $rGoBack.End
$bm2 = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_GoBack", $bm2)
